# feat: Implement exact search on movies name
# Adds a new "Batman" row so an exact-name search picks the most popular
# movie when duplicate names exist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for "Batman"
$ws.Range("A8").Value = "Batman"
$ws.Range("B8").Value = 0.001
$ws.Range("C8").Value = 35000000
$ws.Range("D8").Value = 251409241
$ws.Range("E8").Value = 160160000
$ws.Range("F8").Value = 411569241

# Reflect the (multi-area) selection recorded after adding the row.
$ws.Range("F3").Select()
